# Reorder the priority rows 3-6 on Sheet1:
#  - "README.md" moves from row 3 down to row 6 (with new Importance/Difficulty values)
#  - "Investigate graph issue: redundant append" moves from row 4 up to row 3
#  - "Instument files refer to external cal source" moves from row 5 up to row 4
#  - "Freeze if point moves during scroll" moves from row 6 up to row 5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 3 <- old row 4 values
$ws.Range("A3").Value = "Investigate graph issue: redundant append"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Formula = "=B3*C3/D3"

# New row 4 <- old row 5 values
$ws.Range("A4").Value = "Instument files refer to external cal source"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 3
$ws.Range("E4").Formula = "=B4*C4/D4"

# New row 5 <- old row 6 values
$ws.Range("A5").Value = "Freeze if point moves during scroll"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 3
$ws.Range("E5").Formula = "=B5*C5/D5"

# New row 6 <- old row 3's item, but re-scored after checking time/Assured values
$ws.Range("A6").Value = "README.md"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Formula = "=B6*C6/D6"

# Update the active selection to reflect where the editor left off
$ws.Range("A4").Select()

$wb.Save()
